# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" worksheet (fund-level holdings) right before
#    the "总计" (total) summary sheet.
# 2) Insert a new summary row for "2022-Q1" at the top of the "总计"
#    sheet's data (row 2), pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: new "2022-Q1" sheet with per-fund holding details
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($total)
$q1.Name = "2022-Q1"

# NOTE: sheet references resolve by position, and inserting a sheet at
# $total's old slot leaves $total pointing at the newly-added sheet
# instead of following the "总计" sheet to its new slot. Re-fetch it by
# name before using it again.
$total = $wb.Worksheets.Item("总计")

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Copy the existing bold/centered/bordered header + index-column look from
# the "总计" sheet (chained property setters proved unreliable here).
$total.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q1.Range("A2:A11").PasteSpecial(-4122)

# Force columns B and D:G to be stored as plain text: fund codes like
# "009937" must keep their leading zero, and values like "2.80" must
# keep their trailing zero, instead of being normalized to a number.
$q1.Range("B2:B11").NumberFormat = "@"
$q1.Range("D2:G11").NumberFormat = "@"

$rows = @(
    @("400001", "东方龙混合",                     "2.80", "84.04", "3.39", "0.0949", 10),
    @("009937", "东方欣益一年持有期偏债混合A",       "3.39", "31.18", "1.78", "0.0603", 5),
    @("014125", "华夏中证1000指数增强A",            "7.03", "89.75", "0.82", "0.0576", 7),
    @("014126", "华夏中证1000指数增强C",            "6.09", "89.75", "0.82", "0.0499", 7),
    @("519097", "新华中小市值优选混合",              "0.75", "62.70", "5.58", "0.0418", 1),
    @("006123", "中融高股息精选混合A",               "0.58", "92.22", "2.25", "0.0130", 10),
    @("009938", "东方欣益一年持有期偏债混合C",        "0.52", "31.18", "1.78", "0.0093", 5),
    @("006124", "中融高股息精选混合C",               "0.23", "92.22", "2.25", "0.0052", 10),
    @("001273", "民生加银新动力灵活配置混合A",        "0.04", "68.44", "1.54", "0.0006", 10),
    @("001274", "民生加银新动力灵活配置混合D",        "0.04", "68.44", "1.54", "0.0006", 10)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $q1.Cells.Item($r, 1).Value = $i
    $q1.Cells.Item($r, 2).Value = $data[0]
    $q1.Cells.Item($r, 3).Value = $data[1]
    $q1.Cells.Item($r, 4).Value = $data[2]
    $q1.Cells.Item($r, 5).Value = $data[3]
    $q1.Cells.Item($r, 6).Value = $data[4]
    $q1.Cells.Item($r, 7).Value = $data[5]
    $q1.Cells.Item($r, 8).Value = $data[6]
}

# ---------------------------------------------------------------------
# Part 2: add a "2022-Q1" row to the "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").Style = "Normal"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 10
$total.Range("D2").Value = 0.33
